$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 20 new data rows (rows 323-342) to the dataset, continuing the
# "2-chloroalkyl-tetrahydroxy-p-tert-butylcalix[4]arene" extractant series.

$r = 323
$ws.Cells.Item($r, 1).Value = 321
$ws.Cells.Item($r, 2).Value = "2-chloroalkyl-tetrahydroxy-p-tert-butylcalix[4]arene"
$ws.Cells.Item($r, 3).Value = 4
$ws.Cells.Item($r, 4).Value = "t-butyl"
$ws.Cells.Item($r, 5).Value = "H"
$ws.Cells.Item($r, 6).Value = "(CH2)6Cl"
$ws.Cells.Item($r, 7).Value = "Gd(III)"
$ws.Cells.Item($r, 8).Value = 0.001425
$ws.Cells.Item($r, 9).Value = 0.0000057
$ws.Cells.Item($r, 10).Value = "NO3(2-), dichloromethane, H2O, citric acid, Na2HPO8"
$ws.Cells.Item($r, 11).Value = "Na2HPO4 is disodium phosphate; (CH2)6Cl is the alkyl chloride methylene bridge, which is a linear chain of CH2 groups with a Cl bound at the end"
$ws.Cells.Item($r, 12).Value = "84.15% 0.1 M citric acid/15.85% Na2HPO4 by volume"
$ws.Cells.Item($r, 13).Value = 13.79

$r = 324
$ws.Cells.Item($r, 1).Value = 322
$ws.Cells.Item($r, 2).Value = "2-chloroalkyl-tetrahydroxy-p-tert-butylcalix[4]arene"
$ws.Cells.Item($r, 3).Value = 4
$ws.Cells.Item($r, 4).Value = "t-butyl"
$ws.Cells.Item($r, 5).Value = "H"
$ws.Cells.Item($r, 6).Value = "(CH2)6Cl"
$ws.Cells.Item($r, 7).Value = "Er(III)"
$ws.Cells.Item($r, 8).Value = 0.001425
$ws.Cells.Item($r, 9).Value = 0.0000057
$ws.Cells.Item($r, 10).Value = "NO3(2-), dichloromethane, H2O, citric acid, Na2HPO8"
$ws.Cells.Item($r, 11).Value = "Na2HPO4 is disodium phosphate; (CH2)6Cl is the alkyl chloride methylene bridge, which is a linear chain of CH2 groups with a Cl bound at the end"
$ws.Cells.Item($r, 12).Value = "84.15% 0.1 M citric acid/15.85% Na2HPO4 by volume"
$ws.Cells.Item($r, 13).Value = 11.11

$r = 325
$ws.Cells.Item($r, 1).Value = 323
$ws.Cells.Item($r, 2).Value = "2-chloroalkyl-tetrahydroxy-p-tert-butylcalix[4]arene"
$ws.Cells.Item($r, 3).Value = 4
$ws.Cells.Item($r, 4).Value = "t-butyl"
$ws.Cells.Item($r, 5).Value = "H"
$ws.Cells.Item($r, 6).Value = "(CH2)6Cl"
$ws.Cells.Item($r, 7).Value = "Y(III)"
$ws.Cells.Item($r, 8).Value = 0.001425
$ws.Cells.Item($r, 9).Value = 0.0000057
$ws.Cells.Item($r, 10).Value = "NO3(2-), dichloromethane, H2O, citric acid, Na2HPO8"
$ws.Cells.Item($r, 11).Value = "Na2HPO4 is disodium phosphate; (CH2)6Cl is the alkyl chloride methylene bridge, which is a linear chain of CH2 groups with a Cl bound at the end"
$ws.Cells.Item($r, 12).Value = "84.15% 0.1 M citric acid/15.85% Na2HPO4 by volume"
$ws.Cells.Item($r, 13).Value = 7.69

$r = 326
$ws.Cells.Item($r, 1).Value = 324
$ws.Cells.Item($r, 2).Value = "2-chloroalkyl-tetrahydroxy-p-tert-butylcalix[4]arene"
$ws.Cells.Item($r, 3).Value = 4
$ws.Cells.Item($r, 4).Value = "t-butyl"
$ws.Cells.Item($r, 5).Value = "H"
$ws.Cells.Item($r, 6).Value = "(CH2)6Cl"
$ws.Cells.Item($r, 7).Value = "Ho(III)"
$ws.Cells.Item($r, 8).Value = 0.001425
$ws.Cells.Item($r, 9).Value = 0.0000057
$ws.Cells.Item($r, 10).Value = "NO3(2-), dichloromethane, H2O, citric acid, Na2HPO8"
$ws.Cells.Item($r, 11).Value = "Na2HPO4 is disodium phosphate; (CH2)6Cl is the alkyl chloride methylene bridge, which is a linear chain of CH2 groups with a Cl bound at the end"
$ws.Cells.Item($r, 12).Value = "84.15% 0.1 M citric acid/15.85% Na2HPO4 by volume"
$ws.Cells.Item($r, 13).Value = 12.5

$r = 327
$ws.Cells.Item($r, 1).Value = 325
$ws.Cells.Item($r, 2).Value = "2-chloroalkyl-tetrahydroxy-p-tert-butylcalix[4]arene"
$ws.Cells.Item($r, 3).Value = 4
$ws.Cells.Item($r, 4).Value = "t-butyl"
$ws.Cells.Item($r, 5).Value = "H"
$ws.Cells.Item($r, 6).Value = "(CH2)6Cl"
$ws.Cells.Item($r, 7).Value = "Lu(III)"
$ws.Cells.Item($r, 8).Value = 0.001425
$ws.Cells.Item($r, 9).Value = 0.0000057
$ws.Cells.Item($r, 10).Value = "NO3(2-), dichloromethane, H2O, citric acid, Na2HPO8"
$ws.Cells.Item($r, 11).Value = "Na2HPO4 is disodium phosphate; (CH2)6Cl is the alkyl chloride methylene bridge, which is a linear chain of CH2 groups with a Cl bound at the end"
$ws.Cells.Item($r, 12).Value = "84.15% 0.1 M citric acid/15.85% Na2HPO4 by volume"
$ws.Cells.Item($r, 13).Value = 0

$r = 328
$ws.Cells.Item($r, 1).Value = 326
$ws.Cells.Item($r, 2).Value = "2-chloroalkyl-tetrahydroxy-p-tert-butylcalix[4]arene"
$ws.Cells.Item($r, 3).Value = 4
$ws.Cells.Item($r, 4).Value = "t-butyl"
$ws.Cells.Item($r, 5).Value = "H"
$ws.Cells.Item($r, 6).Value = "(CH2)6Cl"
$ws.Cells.Item($r, 7).Value = "Dy(III)"
$ws.Cells.Item($r, 8).Value = 0.001425
$ws.Cells.Item($r, 9).Value = 0.0000057
$ws.Cells.Item($r, 10).Value = "NO3(2-), dichloromethane, H2O, citric acid, Na2HPO8"
$ws.Cells.Item($r, 11).Value = "Na2HPO4 is disodium phosphate; (CH2)6Cl is the alkyl chloride methylene bridge, which is a linear chain of CH2 groups with a Cl bound at the end"
$ws.Cells.Item($r, 12).Value = "84.15% 0.1 M citric acid/15.85% Na2HPO4 by volume"
$ws.Cells.Item($r, 13).Value = 0

$r = 329
$ws.Cells.Item($r, 1).Value = 327
$ws.Cells.Item($r, 2).Value = "2-chloroalkyl-tetrahydroxy-p-tert-butylcalix[4]arene"
$ws.Cells.Item($r, 3).Value = 4
$ws.Cells.Item($r, 4).Value = "t-butyl"
$ws.Cells.Item($r, 5).Value = "H"
$ws.Cells.Item($r, 6).Value = "(CH2)6Cl"
$ws.Cells.Item($r, 7).Value = "Yb(III)"
$ws.Cells.Item($r, 8).Value = 0.001425
$ws.Cells.Item($r, 9).Value = 0.0000057
$ws.Cells.Item($r, 10).Value = "NO3(2-), dichloromethane, H2O, citric acid, Na2HPO8"
$ws.Cells.Item($r, 11).Value = "Na2HPO4 is disodium phosphate; (CH2)6Cl is the alkyl chloride methylene bridge, which is a linear chain of CH2 groups with a Cl bound at the end"
$ws.Cells.Item($r, 12).Value = "84.15% 0.1 M citric acid/15.85% Na2HPO4 by volume"
$ws.Cells.Item($r, 13).Value = 9.09

$r = 330
$ws.Cells.Item($r, 1).Value = 328
$ws.Cells.Item($r, 2).Value = "2-chloroalkyl-tetrahydroxy-p-tert-butylcalix[4]arene"
$ws.Cells.Item($r, 3).Value = 4
$ws.Cells.Item($r, 4).Value = "t-butyl"
$ws.Cells.Item($r, 5).Value = "OCH3"
$ws.Cells.Item($r, 6).Value = "(CH2)6Cl"
$ws.Cells.Item($r, 7).Value = "La(III)"
$ws.Cells.Item($r, 8).Value = 0.001425
$ws.Cells.Item($r, 9).Value = 0.0000057
$ws.Cells.Item($r, 10).Value = "NO3(2-), dichloromethane, H2O, citric acid, Na2HPO4"
$ws.Cells.Item($r, 11).Value = "Na2HPO4 is disodium phosphate; (CH2)6Cl is the alkyl chloride methylene bridge, which is a linear chain of CH2 groups with a Cl bound at the end"
$ws.Cells.Item($r, 12).Value = "84.15% 0.1 M citric acid/15.85% Na2HPO4 by volume"
$ws.Cells.Item($r, 13).Value = 50

$r = 331
$ws.Cells.Item($r, 1).Value = 329
$ws.Cells.Item($r, 2).Value = "2-chloroalkyl-tetrahydroxy-p-tert-butylcalix[4]arene"
$ws.Cells.Item($r, 3).Value = 4
$ws.Cells.Item($r, 4).Value = "t-butyl"
$ws.Cells.Item($r, 5).Value = "OCH3"
$ws.Cells.Item($r, 6).Value = "(CH2)6Cl"
$ws.Cells.Item($r, 7).Value = "Ce(III)"
$ws.Cells.Item($r, 8).Value = 0.001425
$ws.Cells.Item($r, 9).Value = 0.0000057
$ws.Cells.Item($r, 10).Value = "NO3(2-), dichloromethane, H2O, citric acid, Na2HPO4"
$ws.Cells.Item($r, 11).Value = "Na2HPO4 is disodium phosphate; (CH2)6Cl is the alkyl chloride methylene bridge, which is a linear chain of CH2 groups with a Cl bound at the end"
$ws.Cells.Item($r, 12).Value = "84.15% 0.1 M citric acid/15.85% Na2HPO4 by volume"
$ws.Cells.Item($r, 13).Value = 75

$r = 332
$ws.Cells.Item($r, 1).Value = 330
$ws.Cells.Item($r, 2).Value = "2-chloroalkyl-tetrahydroxy-p-tert-butylcalix[4]arene"
$ws.Cells.Item($r, 3).Value = 4
$ws.Cells.Item($r, 4).Value = "t-butyl"
$ws.Cells.Item($r, 5).Value = "OCH3"
$ws.Cells.Item($r, 6).Value = "(CH2)6Cl"
$ws.Cells.Item($r, 7).Value = "Pr(III)"
$ws.Cells.Item($r, 8).Value = 0.001425
$ws.Cells.Item($r, 9).Value = 0.0000057
$ws.Cells.Item($r, 10).Value = "NO3(2-), dichloromethane, H2O, citric acid, Na2HPO5"
$ws.Cells.Item($r, 11).Value = "Na2HPO4 is disodium phosphate; (CH2)6Cl is the alkyl chloride methylene bridge, which is a linear chain of CH2 groups with a Cl bound at the end"
$ws.Cells.Item($r, 12).Value = "84.15% 0.1 M citric acid/15.85% Na2HPO4 by volume"
$ws.Cells.Item($r, 13).Value = 36.36

$r = 333
$ws.Cells.Item($r, 1).Value = 331
$ws.Cells.Item($r, 2).Value = "2-chloroalkyl-tetrahydroxy-p-tert-butylcalix[4]arene"
$ws.Cells.Item($r, 3).Value = 4
$ws.Cells.Item($r, 4).Value = "t-butyl"
$ws.Cells.Item($r, 5).Value = "OCH3"
$ws.Cells.Item($r, 6).Value = "(CH2)6Cl"
$ws.Cells.Item($r, 7).Value = "Nd(III)"
$ws.Cells.Item($r, 8).Value = 0.001425
$ws.Cells.Item($r, 9).Value = 0.0000057
$ws.Cells.Item($r, 10).Value = "NO3(2-), dichloromethane, H2O, citric acid, Na2HPO6"
$ws.Cells.Item($r, 11).Value = "Na2HPO4 is disodium phosphate; (CH2)6Cl is the alkyl chloride methylene bridge, which is a linear chain of CH2 groups with a Cl bound at the end"
$ws.Cells.Item($r, 12).Value = "84.15% 0.1 M citric acid/15.85% Na2HPO4 by volume"
$ws.Cells.Item($r, 13).Value = 38.46

$r = 334
$ws.Cells.Item($r, 1).Value = 332
$ws.Cells.Item($r, 2).Value = "2-chloroalkyl-tetrahydroxy-p-tert-butylcalix[4]arene"
$ws.Cells.Item($r, 3).Value = 4
$ws.Cells.Item($r, 4).Value = "t-butyl"
$ws.Cells.Item($r, 5).Value = "OCH3"
$ws.Cells.Item($r, 6).Value = "(CH2)6Cl"
$ws.Cells.Item($r, 7).Value = "Sm(III)"
$ws.Cells.Item($r, 8).Value = 0.001425
$ws.Cells.Item($r, 9).Value = 0.0000057
$ws.Cells.Item($r, 10).Value = "NO3(2-), dichloromethane, H2O, citric acid, Na2HPO7"
$ws.Cells.Item($r, 11).Value = "Na2HPO4 is disodium phosphate; (CH2)6Cl is the alkyl chloride methylene bridge, which is a linear chain of CH2 groups with a Cl bound at the end"
$ws.Cells.Item($r, 12).Value = "84.15% 0.1 M citric acid/15.85% Na2HPO4 by volume"
$ws.Cells.Item($r, 13).Value = 7.69

$r = 335
$ws.Cells.Item($r, 1).Value = 333
$ws.Cells.Item($r, 2).Value = "2-chloroalkyl-tetrahydroxy-p-tert-butylcalix[4]arene"
$ws.Cells.Item($r, 3).Value = 4
$ws.Cells.Item($r, 4).Value = "t-butyl"
$ws.Cells.Item($r, 5).Value = "OCH3"
$ws.Cells.Item($r, 6).Value = "(CH2)6Cl"
$ws.Cells.Item($r, 7).Value = "Eu(III)"
$ws.Cells.Item($r, 8).Value = 0.001425
$ws.Cells.Item($r, 9).Value = 0.0000057
$ws.Cells.Item($r, 10).Value = "NO3(2-), dichloromethane, H2O, citric acid, Na2HPO8"
$ws.Cells.Item($r, 11).Value = "Na2HPO4 is disodium phosphate; (CH2)6Cl is the alkyl chloride methylene bridge, which is a linear chain of CH2 groups with a Cl bound at the end"
$ws.Cells.Item($r, 12).Value = "84.15% 0.1 M citric acid/15.85% Na2HPO4 by volume"
$ws.Cells.Item($r, 13).Value = 27.27

$r = 336
$ws.Cells.Item($r, 1).Value = 334
$ws.Cells.Item($r, 2).Value = "2-chloroalkyl-tetrahydroxy-p-tert-butylcalix[4]arene"
$ws.Cells.Item($r, 3).Value = 4
$ws.Cells.Item($r, 4).Value = "t-butyl"
$ws.Cells.Item($r, 5).Value = "OCH3"
$ws.Cells.Item($r, 6).Value = "(CH2)6Cl"
$ws.Cells.Item($r, 7).Value = "Gd(III)"
$ws.Cells.Item($r, 8).Value = 0.001425
$ws.Cells.Item($r, 9).Value = 0.0000057
$ws.Cells.Item($r, 10).Value = "NO3(2-), dichloromethane, H2O, citric acid, Na2HPO8"
$ws.Cells.Item($r, 11).Value = "Na2HPO4 is disodium phosphate; (CH2)6Cl is the alkyl chloride methylene bridge, which is a linear chain of CH2 groups with a Cl bound at the end"
$ws.Cells.Item($r, 12).Value = "84.15% 0.1 M citric acid/15.85% Na2HPO4 by volume"
$ws.Cells.Item($r, 13).Value = 0

$r = 337
$ws.Cells.Item($r, 1).Value = 335
$ws.Cells.Item($r, 2).Value = "2-chloroalkyl-tetrahydroxy-p-tert-butylcalix[4]arene"
$ws.Cells.Item($r, 3).Value = 4
$ws.Cells.Item($r, 4).Value = "t-butyl"
$ws.Cells.Item($r, 5).Value = "OCH3"
$ws.Cells.Item($r, 6).Value = "(CH2)6Cl"
$ws.Cells.Item($r, 7).Value = "Er(III)"
$ws.Cells.Item($r, 8).Value = 0.001425
$ws.Cells.Item($r, 9).Value = 0.0000057
$ws.Cells.Item($r, 10).Value = "NO3(2-), dichloromethane, H2O, citric acid, Na2HPO8"
$ws.Cells.Item($r, 11).Value = "Na2HPO4 is disodium phosphate; (CH2)6Cl is the alkyl chloride methylene bridge, which is a linear chain of CH2 groups with a Cl bound at the end"
$ws.Cells.Item($r, 12).Value = "84.15% 0.1 M citric acid/15.85% Na2HPO4 by volume"
$ws.Cells.Item($r, 13).Value = 0

$r = 338
$ws.Cells.Item($r, 1).Value = 336
$ws.Cells.Item($r, 2).Value = "2-chloroalkyl-tetrahydroxy-p-tert-butylcalix[4]arene"
$ws.Cells.Item($r, 3).Value = 4
$ws.Cells.Item($r, 4).Value = "t-butyl"
$ws.Cells.Item($r, 5).Value = "OCH3"
$ws.Cells.Item($r, 6).Value = "(CH2)6Cl"
$ws.Cells.Item($r, 7).Value = "Y(III)"
$ws.Cells.Item($r, 8).Value = 0.001425
$ws.Cells.Item($r, 9).Value = 0.0000057
$ws.Cells.Item($r, 10).Value = "NO3(2-), dichloromethane, H2O, citric acid, Na2HPO8"
$ws.Cells.Item($r, 11).Value = "Na2HPO4 is disodium phosphate; (CH2)6Cl is the alkyl chloride methylene bridge, which is a linear chain of CH2 groups with a Cl bound at the end"
$ws.Cells.Item($r, 12).Value = "84.15% 0.1 M citric acid/15.85% Na2HPO4 by volume"
$ws.Cells.Item($r, 13).Value = 7.14

$r = 339
$ws.Cells.Item($r, 1).Value = 337
$ws.Cells.Item($r, 2).Value = "2-chloroalkyl-tetrahydroxy-p-tert-butylcalix[4]arene"
$ws.Cells.Item($r, 3).Value = 4
$ws.Cells.Item($r, 4).Value = "t-butyl"
$ws.Cells.Item($r, 5).Value = "OCH3"
$ws.Cells.Item($r, 6).Value = "(CH2)6Cl"
$ws.Cells.Item($r, 7).Value = "Ho(III)"
$ws.Cells.Item($r, 8).Value = 0.001425
$ws.Cells.Item($r, 9).Value = 0.0000057
$ws.Cells.Item($r, 10).Value = "NO3(2-), dichloromethane, H2O, citric acid, Na2HPO8"
$ws.Cells.Item($r, 11).Value = "Na2HPO4 is disodium phosphate; (CH2)6Cl is the alkyl chloride methylene bridge, which is a linear chain of CH2 groups with a Cl bound at the end"
$ws.Cells.Item($r, 12).Value = "84.15% 0.1 M citric acid/15.85% Na2HPO4 by volume"
$ws.Cells.Item($r, 13).Value = 15.38

$r = 340
$ws.Cells.Item($r, 1).Value = 338
$ws.Cells.Item($r, 2).Value = "2-chloroalkyl-tetrahydroxy-p-tert-butylcalix[4]arene"
$ws.Cells.Item($r, 3).Value = 4
$ws.Cells.Item($r, 4).Value = "t-butyl"
$ws.Cells.Item($r, 5).Value = "OCH3"
$ws.Cells.Item($r, 6).Value = "(CH2)6Cl"
$ws.Cells.Item($r, 7).Value = "Lu(III)"
$ws.Cells.Item($r, 8).Value = 0.001425
$ws.Cells.Item($r, 9).Value = 0.0000057
$ws.Cells.Item($r, 10).Value = "NO3(2-), dichloromethane, H2O, citric acid, Na2HPO8"
$ws.Cells.Item($r, 11).Value = "Na2HPO4 is disodium phosphate; (CH2)6Cl is the alkyl chloride methylene bridge, which is a linear chain of CH2 groups with a Cl bound at the end"
$ws.Cells.Item($r, 12).Value = "84.15% 0.1 M citric acid/15.85% Na2HPO4 by volume"
$ws.Cells.Item($r, 13).Value = 0

$r = 341
$ws.Cells.Item($r, 1).Value = 339
$ws.Cells.Item($r, 2).Value = "2-chloroalkyl-tetrahydroxy-p-tert-butylcalix[4]arene"
$ws.Cells.Item($r, 3).Value = 4
$ws.Cells.Item($r, 4).Value = "t-butyl"
$ws.Cells.Item($r, 5).Value = "OCH3"
$ws.Cells.Item($r, 6).Value = "(CH2)6Cl"
$ws.Cells.Item($r, 7).Value = "Dy(III)"
$ws.Cells.Item($r, 8).Value = 0.001425
$ws.Cells.Item($r, 9).Value = 0.0000057
$ws.Cells.Item($r, 10).Value = "NO3(2-), dichloromethane, H2O, citric acid, Na2HPO8"
$ws.Cells.Item($r, 11).Value = "Na2HPO4 is disodium phosphate; (CH2)6Cl is the alkyl chloride methylene bridge, which is a linear chain of CH2 groups with a Cl bound at the end"
$ws.Cells.Item($r, 12).Value = "84.15% 0.1 M citric acid/15.85% Na2HPO4 by volume"
$ws.Cells.Item($r, 13).Value = 15

$r = 342
$ws.Cells.Item($r, 1).Value = 340
$ws.Cells.Item($r, 2).Value = "2-chloroalkyl-tetrahydroxy-p-tert-butylcalix[4]arene"
$ws.Cells.Item($r, 3).Value = 4
$ws.Cells.Item($r, 4).Value = "t-butyl"
$ws.Cells.Item($r, 5).Value = "OCH3"
$ws.Cells.Item($r, 6).Value = "(CH2)6Cl"
$ws.Cells.Item($r, 7).Value = "Yb(III)"
$ws.Cells.Item($r, 8).Value = 0.001425
$ws.Cells.Item($r, 9).Value = 0.0000057
$ws.Cells.Item($r, 10).Value = "NO3(2-), dichloromethane, H2O, citric acid, Na2HPO8"
$ws.Cells.Item($r, 11).Value = "Na2HPO4 is disodium phosphate; (CH2)6Cl is the alkyl chloride methylene bridge, which is a linear chain of CH2 groups with a Cl bound at the end"
$ws.Cells.Item($r, 12).Value = "84.15% 0.1 M citric acid/15.85% Na2HPO4 by volume"
$ws.Cells.Item($r, 13).Value = 9.09

# Update the active selection / scrolled position to match where the
# user ended up after entering the new data.
$ws.Range("A322").Select()
